# Apply the dropdown-code-generator edits to Sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the I2 label and convert the numeric I-column sample values into
# text values (order matches the shared-string table growth: hi, hello,
# rawr, v1, v2).
$ws.Range("I2").Value = "hi"
$ws.Range("I3").Value = "hello"
$ws.Range("I4").Value = "rawr"
$ws.Range("I5").Value = "hi"

# Update the J2:K2 header/label cells to new text values.
$ws.Range("J2").Value = "v1"
$ws.Range("K2").Value = "v2"

# Update the active selection on the sheet to reflect the new working cell.
$ws.Activate()
$ws.Range("J6").Select()
